$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $rng = $ws.Range($cell)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

Set-TextValue "D2" "292.36"
Set-TextValue "E2" "-6.97%"
Set-TextValue "D3" "40.42"
Set-TextValue "E3" "0.39%"
Set-TextValue "D4" "5.019"
Set-TextValue "E4" "-2.88%"
Set-TextValue "D5" "0.07319"
Set-TextValue "B6" "GateToken"
Set-TextValue "C6" "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
Set-TextValue "D6" "4.300"
Set-TextValue "E6" "-0.62%"
Set-TextValue "B7" "FTXToken"
Set-TextValue "C7" "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
Set-TextValue "D7" "1.531"
Set-TextValue "E7" "-7.99%"
Set-TextValue "B8" "MXToken"
Set-TextValue "C8" "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
Set-TextValue "D8" "0.9270"
Set-TextValue "E8" "-0.03%"
Set-TextValue "B9" "BTSEToken"
Set-TextValue "C9" "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
Set-TextValue "D9" "2.369"
Set-TextValue "E9" "-2.27%"
Set-TextValue "B10" "LiechtensteinCryptoassetsExchange"
Set-TextValue "C10" "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
Set-TextValue "D10" "0.1186"
Set-TextValue "E10" "-1.03%"
Set-TextValue "B11" "WazirX"
Set-TextValue "C11" "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
Set-TextValue "D11" "0.1743"
Set-TextValue "E11" "-4.30%"
Set-TextValue "B12" "BitrueCoin"
Set-TextValue "C12" "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
Set-TextValue "D12" "0.04324"
Set-TextValue "E12" "4.51%"
Set-TextValue "B13" "MandalaExchangeToken"
Set-TextValue "C13" "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
Set-TextValue "D13" "0.08695"
Set-TextValue "E13" "-3.61%"
Set-TextValue "B14" "BitMartToken"
Set-TextValue "C14" "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
Set-TextValue "D14" "0.1054"
Set-TextValue "E14" "0.07%"
Set-TextValue "B15" "BitForexToken"
Set-TextValue "C15" "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
Set-TextValue "D15" "0.001265"
Set-TextValue "E15" "-1.67%"
Set-TextValue "B16" "TigerCash"
Set-TextValue "C16" "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
Set-TextValue "D16" "0.006005"
Set-TextValue "E16" "3.54%"
Set-TextValue "B17" "LEO"
Set-TextValue "C17" "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
Set-TextValue "D17" "3.339"
Set-TextValue "E17" "-0.37%"
Set-TextValue "B18" "BitpandaEcosystemToken"
Set-TextValue "C18" "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
Set-TextValue "D18" "0.3289"
Set-TextValue "E18" "-2.02%"
Set-TextValue "B19" "MCDex"
Set-TextValue "C19" "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
Set-TextValue "D19" "7.972"
Set-TextValue "E19" "5.52%"
Set-TextValue "B20" "ProBitToken"
Set-TextValue "C20" "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
Set-TextValue "D20" "0.1391"
Set-TextValue "E20" "2.85%"
Set-TextValue "B21" "ZBToken"
Set-TextValue "C21" "https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb"
Set-TextValue "D21" "0.2793"
Set-TextValue "E21" "-0.42%"
Set-TextValue "D22" "0.03931"
Set-TextValue "E22" "-2.60%"
Set-TextValue "E23" "-1.02%"
Set-TextValue "E24" "-7.13%"
Set-TextValue "E25" "0.88%"
Set-TextValue "D26" "0.0003726"
Set-TextValue "D38" "0.02277"
Set-TextValue "E38" "-5.67%"
Set-TextValue "D39" "0.04977"
Set-TextValue "E39" "-3.27%"
Set-TextValue "E40" "70.55%"
Set-TextValue "D41" "0.007716"
Set-TextValue "E41" "-0.07%"
Set-TextValue "D42" "0.1284"
Set-TextValue "E42" "-1.24%"
Set-TextValue "D43" "0.007366"
Set-TextValue "E43" "-3.07%"
Set-TextValue "D44" "0.007303"
Set-TextValue "E44" "-16.35%"
Set-TextValue "D45" "0.2921"
Set-TextValue "E45" "-13.96%"
Set-TextValue "D46" "0.00006300"
Set-TextValue "E46" "-4.45%"
Set-TextValue "E47" "0.03%"
Set-TextValue "E48" "-92.07%"
Set-TextValue "E49" "0.03%"
Set-TextValue "E50" "0.03%"

Write-Host "Applied all cell updates"